# Rename worksheets: replace the first underscore in each sheet name with a space.
$wb = $excel.ActiveWorkbook

$renames = @{
    "Global_Concordance" = "Global Concordance"
    "Discordance_SJR-index" = "Discordance SJR-index"
    "Discordance_Cite Score" = "Discordance Cite Score"
    "Discordance_H-index" = "Discordance H-index"
    "Discordance_Best Subject Rank" = "Discordance Best Subject Rank"
    "Discordance_Total Docs" = "Discordance Total Docs"
    "Discordance_Total Docs 3y" = "Discordance Total Docs 3y"
    "Discordance_Total Refs" = "Discordance Total Refs"
    "Discordance_Total Cites 3y" = "Discordance Total Cites 3y"
    "Discordance_Citable Docs 3y" = "Discordance Citable Docs 3y"
    "Discordance_Cites per Doc 2y" = "Discordance Cites per Doc 2y"
    "Discordance_Refs per Doc" = "Discordance Refs per Doc"
    "Discordance_Coverage" = "Discordance Coverage"
    "Final_Ranking" = "Final Ranking"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renames.ContainsKey($oldName)) {
        $ws.Name = $renames[$oldName]
    }
}
